$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 4) with the date_col configuration entry
$ws.Range("A4").Value = "date_col"
$ws.Range("B4").Value = "issue_date,application_time"

# Move the active selection to A5, matching the post-edit cursor position
$ws.Range("A5").Select()
